# Auto-generated edit script: updates Leve profit calculation values
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# matching an upstream data refresh from the scheduled price-scrape runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 363
$ws.Range("I2").Value = 224.33333
$ws.Range("K2").Value = 224.33333
$ws.Range("M2").Value = -111.33333
$ws.Range("H9").Value = 267.5
$ws.Range("I9").Value = 267.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 267.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -98.5
$ws.Range("N9").ClearContents()
$ws.Range("H18").Value = 789.5
$ws.Range("I18").Value = 789.5
$ws.Range("K18").Value = 789.5
$ws.Range("M18").Value = -505.5
$ws.Range("H43").Value = 2152
$ws.Range("I43").Value = 1501
$ws.Range("J43").Value = 2477.5
$ws.Range("K43").Value = 1501
$ws.Range("L43").Value = 2477.5
$ws.Range("M43").Value = -1432
$ws.Range("N43").Value = -2615.5
$ws.Range("H115").Value = 1516.25
$ws.Range("I115").Value = 1018.7143
$ws.Range("K115").Value = 3056.1429
$ws.Range("M115").Value = -1489.1429
$ws.Range("H135").Value = 1858.8
$ws.Range("I135").Value = 1490.25
$ws.Range("K135").Value = 13412.25
$ws.Range("M135").Value = -10877.25
$ws.Range("H138").Value = 3099.3333
$ws.Range("I138").Value = 2842.25
$ws.Range("J138").Value = 3613.5
$ws.Range("K138").Value = 8526.75
$ws.Range("L138").Value = 10840.5
$ws.Range("M138").Value = -3386.75
$ws.Range("N138").Value = -21120.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2456
$ws.Range("I32").Value = 2456
$ws.Range("K32").Value = 2456
$ws.Range("M32").Value = -2169
$ws.Range("H61").Value = 5479.8
$ws.Range("I61").Value = 4133
$ws.Range("K61").Value = 4133
$ws.Range("M61").Value = -3921
$ws.Range("H74").Value = 5997
$ws.Range("I74").Value = 1995
$ws.Range("K74").Value = 1995
$ws.Range("M74").Value = -1121
$ws.Range("H77").Value = 5997
$ws.Range("I77").Value = 1995
$ws.Range("K77").Value = 9975
$ws.Range("M77").Value = -5607
$ws.Range("H136").Value = 5479.8
$ws.Range("I136").Value = 4133
$ws.Range("K136").Value = 12399
$ws.Range("M136").Value = -9849

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 22640
$ws.Range("I75").Value = 1233.3334
$ws.Range("K75").Value = 1233.3334
$ws.Range("M75").Value = -297.3334
$ws.Range("H78").Value = 22640
$ws.Range("I78").Value = 1233.3334
$ws.Range("K78").Value = 3700.0002
$ws.Range("M78").Value = 979.9998000000001
$ws.Range("H99").Value = 3033
$ws.Range("I99").Value = 2513.2222
$ws.Range("J99").Value = 5372
$ws.Range("K99").Value = 2513.2222
$ws.Range("L99").Value = 5372
$ws.Range("M99").Value = -1015.2222
$ws.Range("N99").Value = -8368

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5998.6665
$ws.Range("I31").Value = 2618.7856
$ws.Range("J31").Value = 9638.538
$ws.Range("K31").Value = 2618.7856
$ws.Range("L31").Value = 9638.538
$ws.Range("M31").Value = -2323.7856
$ws.Range("N31").Value = -10228.538
$ws.Range("H34").Value = 5998.6665
$ws.Range("I34").Value = 2618.7856
$ws.Range("J34").Value = 9638.538
$ws.Range("K34").Value = 2618.7856
$ws.Range("L34").Value = 9638.538
$ws.Range("M34").Value = -2416.7856
$ws.Range("N34").Value = -10042.538
$ws.Range("H58").Value = 3942.3684
$ws.Range("I58").Value = 2961.1333
$ws.Range("J58").Value = 7622
$ws.Range("K58").Value = 2961.1333
$ws.Range("L58").Value = 7622
$ws.Range("M58").Value = -2758.1333
$ws.Range("N58").Value = -8028
$ws.Range("H86").Value = 1969
$ws.Range("I86").Value = 1953.5
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1953.5
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -830.5
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1969
$ws.Range("I89").Value = 1953.5
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 9767.5
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -4151.5
$ws.Range("N89").Value = -21232
$ws.Range("H105").Value = 3151.25
$ws.Range("I105").Value = 2720.3333
$ws.Range("K105").Value = 2720.3333
$ws.Range("M105").Value = -973.3332999999998
$ws.Range("H107").Value = 565.2
$ws.Range("I107").Value = 559.5714
$ws.Range("J107").Value = 644
$ws.Range("K107").Value = 559.5714
$ws.Range("L107").Value = 644
$ws.Range("M107").Value = 1360.4286
$ws.Range("N107").Value = -4484
$ws.Range("H132").Value = 3195.3333
$ws.Range("I132").Value = 2823.2727
$ws.Range("K132").Value = 8469.8181
$ws.Range("M132").Value = -5939.8181
$ws.Range("H136").Value = 3942.3684
$ws.Range("I136").Value = 2961.1333
$ws.Range("J136").Value = 7622
$ws.Range("K136").Value = 8883.3999
$ws.Range("L136").Value = 22866
$ws.Range("M136").Value = -6333.3999
$ws.Range("N136").Value = -27966

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2699.2222
$ws.Range("I113").Value = 1149.5
$ws.Range("J113").Value = 3142
$ws.Range("K113").Value = 3448.5
$ws.Range("L113").Value = 9426
$ws.Range("M113").Value = -1278.5
$ws.Range("N113").Value = -13766
$ws.Range("H129").Value = 2038.9166
$ws.Range("I129").Value = 706.5
$ws.Range("K129").Value = 2119.5
$ws.Range("M129").Value = 2880.5
$ws.Range("H138").Value = 6007.615
$ws.Range("I138").Value = 2585.7144
$ws.Range("J138").Value = 9999.833000000001
$ws.Range("K138").Value = 7757.1432
$ws.Range("L138").Value = 29999.499
$ws.Range("M138").Value = -2617.1432
$ws.Range("N138").Value = -40279.499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 89000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 89000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 89000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -92744
$ws.Range("H122").Value = 719243.7
$ws.Range("I122").Value = 1004393.8
$ws.Range("K122").Value = 3013181.4
$ws.Range("M122").Value = -3010731.4
$ws.Range("H132").Value = 19975.143
$ws.Range("I132").Value = 21654.334
$ws.Range("J132").Value = 9900
$ws.Range("K132").Value = 64963.00199999999
$ws.Range("L132").Value = 29700
$ws.Range("M132").Value = -62433.00199999999
$ws.Range("N132").Value = -34760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H50").Value = 67084
$ws.Range("J50").Value = 67084
$ws.Range("L50").Value = 67084
$ws.Range("N50").Value = -68346
$ws.Range("H132").Value = 4066.4707
$ws.Range("I132").Value = 3820.6875
$ws.Range("J132").Value = 7999
$ws.Range("K132").Value = 11462.0625
$ws.Range("L132").Value = 23997
$ws.Range("M132").Value = -8932.0625
$ws.Range("N132").Value = -29057
